# Auto-generated edit script: updates cryptos list Price (D) and Volume(1h) (E) columns
# per commit "Updated cryptos list on Sun Oct 15 15:15:41 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to remain a text string even when the value looks like
    # a number (e.g. "210.13"), then restore the default/normal style so no
    # stray number-format style attribute gets attached to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.043.90"
Set-TextValue $ws.Range("E2") "  +0.37%  "
Set-TextValue $ws.Range("D3") "1.561.66"
Set-TextValue $ws.Range("E3") "  +0.56%  "
Set-TextValue $ws.Range("E4") "  +0.36%  "
Set-TextValue $ws.Range("D5") "210.13"
Set-TextValue $ws.Range("E5") "  +1.69%  "
Set-TextValue $ws.Range("D6") "0.489"
Set-TextValue $ws.Range("E6") "  +0.29%  "
Set-TextValue $ws.Range("E7") "  +0.46%  "
Set-TextValue $ws.Range("D8") "21.86"
Set-TextValue $ws.Range("E8") "  -0.68%  "
Set-TextValue $ws.Range("E9") "  +0.01%  "
Set-TextValue $ws.Range("D10") "0.0596"
Set-TextValue $ws.Range("E10") "  +0.20%  "
Set-TextValue $ws.Range("E11") "  +0.50%  "
Set-TextValue $ws.Range("D12") "1.782.58"
Set-TextValue $ws.Range("E12") "  +0.55%  "
Set-TextValue $ws.Range("D13") "1.442.97"
Set-TextValue $ws.Range("E13") "  -6.86%  "
Set-TextValue $ws.Range("D14") "3.76"
Set-TextValue $ws.Range("E14") "  +0.36%  "
Set-TextValue $ws.Range("D15") "0.517"
Set-TextValue $ws.Range("E15") "  -0.42%  "
Set-TextValue $ws.Range("D16") "27.045.21"
Set-TextValue $ws.Range("E16") "  +0.48%  "
Set-TextValue $ws.Range("D17") "61.90"
Set-TextValue $ws.Range("E17") "  +0.44%  "
Set-TextValue $ws.Range("D18") "0.0₃0701"
Set-TextValue $ws.Range("E18") "  -0.66%  "
Set-TextValue $ws.Range("D19") "214.81"
Set-TextValue $ws.Range("E19") "  -1.10%  "
Set-TextValue $ws.Range("D20") "7.34"
Set-TextValue $ws.Range("E20") "  +0.62%  "
Set-TextValue $ws.Range("E21") "  +0.48%  "
Set-TextValue $ws.Range("D22") "4.12"
Set-TextValue $ws.Range("E22") "  +0.77%  "
Set-TextValue $ws.Range("D23") "9.16"
Set-TextValue $ws.Range("E23") "  -0.28%  "
Set-TextValue $ws.Range("D24") "1.93"
Set-TextValue $ws.Range("E24") "  -0.47%  "
Set-TextValue $ws.Range("D25") "153.82"
Set-TextValue $ws.Range("E25") "  +0.22%  "
Set-TextValue $ws.Range("D26") "6.58"
Set-TextValue $ws.Range("E26") "  -0.76%  "
Set-TextValue $ws.Range("D27") "15.00"
Set-TextValue $ws.Range("E27") "  +0.17%  "
Set-TextValue $ws.Range("E28") "  +1.15%  "
Set-TextValue $ws.Range("E29") "  +0.46%  "
Set-TextValue $ws.Range("E30") "  +4.45%  "
Set-TextValue $ws.Range("D31") "0.0470"
Set-TextValue $ws.Range("E31") "  +0.23%  "
Set-TextValue $ws.Range("E32") "  +0.39%  "
Set-TextValue $ws.Range("D33") "3.17"
Set-TextValue $ws.Range("E33") "  +1.98%  "
Set-TextValue $ws.Range("D34") "1.425.55"
Set-TextValue $ws.Range("E34") "  +0.52%  "
Set-TextValue $ws.Range("E35") "  +0.42%  "
Set-TextValue $ws.Range("D36") "1.59"
Set-TextValue $ws.Range("E36") "  -0.50%  "
Set-TextValue $ws.Range("D37") "2.34"
Set-TextValue $ws.Range("E37") "  +1.76%  "
Set-TextValue $ws.Range("D38") "0.0166"
Set-TextValue $ws.Range("E38") "  +0.70%  "
Set-TextValue $ws.Range("D39") "0.528"
Set-TextValue $ws.Range("E39") "  +0.38%  "
Set-TextValue $ws.Range("D40") "5.78"
Set-TextValue $ws.Range("E40") "  +2.51%  "
Set-TextValue $ws.Range("D41") "0.804"
Set-TextValue $ws.Range("E41") "  -0.21%  "
Set-TextValue $ws.Range("E42") "  +0.43%  "
Set-TextValue $ws.Range("D43") "2.34"
Set-TextValue $ws.Range("E43") "  +1.47%  "
Set-TextValue $ws.Range("D44") "0.999"
Set-TextValue $ws.Range("E44") "  +0.21%  "
Set-TextValue $ws.Range("D45") "64.28"
Set-TextValue $ws.Range("E45") "  -0.25%  "
Set-TextValue $ws.Range("D46") "1.73"
Set-TextValue $ws.Range("E46") "  +0.04%  "
Set-TextValue $ws.Range("D47") "1.700.98"
Set-TextValue $ws.Range("E47") "  +0.88%  "
Set-TextValue $ws.Range("D48") "85.89"
Set-TextValue $ws.Range("E48") "  -1.51%  "
Set-TextValue $ws.Range("E49") "  +2.30%  "
Set-TextValue $ws.Range("D50") "0.0517"
Set-TextValue $ws.Range("E50") "  -0.50%  "
Set-TextValue $ws.Range("D51") "0.0955"
Set-TextValue $ws.Range("E51") "  -0.31%  "
